$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)          # "Subtitle 4" placeholder on slide 1
$tr = $shp.TextFrame.TextRange
$para = $tr.Paragraphs(4)         # the ".... ~  Amaydeep minz" line

# The paragraph is currently three runs followed by an unaffected fourth run:
#   run1 (chars 1-54):  "." + alternating nbsp/space filler + "~" + nbsp + " "
#   run2 (chars 55-62): "Amaydeep"   (flagged err="1" by the spell checker)
#   run3 (char  63):    " "
#   run4 (chars 64-67): "minz"       <- left untouched
#
# Replacing the text that spans runs 1-3 with the corrected name merges
# those three runs into a single run (keeping run1's formatting, i.e. no
# err="1"), which is exactly what the target edit does: "Amaydeep" becomes
# "Amay deep".

$nbsp = [char]0x00A0
$filler = ""
for ($i = 0; $i -lt 25; $i++) {
    $filler = $filler + $nbsp + " "
}
$newText = "." + $filler + "~" + $nbsp + " " + "Amay deep "

$combined = $para.Characters(1, 63)
$combined.Text = $newText

$para.Text
